$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.34288292259397146
$ws.Range("D2").Value = 0.15
$ws.Range("E2").Value = 0.2
$ws.Range("F2").Value = 2.1
